$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.796.21'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.811.04'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.30'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.83'
$ws.Range('E6').Value = '  +1.34%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.160'
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.06'
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.447.00'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.844.87'
$ws.Range('E15').Value = '  +2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.58'
$ws.Range('E16').Value = '  +5.18%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.779.90'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.10'
$ws.Range('E18').Value = '  +2.71%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '462.22'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.97'
$ws.Range('E21').Value = '  -5.22%  '
$ws.Range('E22').Value = '  +1.27%  '
$ws.Range('E23').Value = '  +2.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.61'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.12'
$ws.Range('E25').Value = '  +2.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.12'
$ws.Range('E26').Value = '  -0.91%  '
$ws.Range('E27').Value = '  +0.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.957.94'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.78'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  +4.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.30'
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.75'
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('B34').Value = 'Aptos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.11'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.749.34'
$ws.Range('E36').Value = '  +0.59%  '
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.44'
$ws.Range('E38').Value = '  +2.37%  '
$ws.Range('E39').Value = '  +0.85%  '
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('E41').Value = '  +1.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '48.19'
$ws.Range('E44').Value = '  +2.94%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.68'
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('E46').Value = '  +0.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '149.30'
$ws.Range('E47').Value = '  +2.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.34'
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '398.06'
$ws.Range('E49').Value = '  +2.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.85'
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.76'
$ws.Range('E51').Value = '  +6.02%  '
